# Change the table style used by the table on slide 16 from the
# presentation's custom "Table_0" style to the built-in PowerPoint table
# style {CEB7269D-37F7-46FF-AC3C-E26E7C18039F}.
#
# <a:tableStyleId>{5BF5068B-CE96-40BA-9859-0B60F4375200}</a:tableStyleId>
#   -> <a:tableStyleId>{CEB7269D-37F7-46FF-AC3C-E26E7C18039F}</a:tableStyleId>

$p = $ppt.ActivePresentation

$targetStyleId = "{CEB7269D-37F7-46FF-AC3C-E26E7C18039F}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($targetStyleId)
        }
    }
}
